$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G1").Value = "singleDevice"
$ws.Range("H1").Value = "deviceId"

$ws.Range("H1").Select()
